# Add a new drug record (row 2) to the "drugs" sheet, then create a new
# "Sheet1" worksheet after it holding the same record transposed into a
# single column (as if the row had been copied and pasted-special/transposed).

$wb = $excel.ActiveWorkbook
$drugs = $wb.Worksheets.Item("drugs")

# --- New drug row on "drugs" ------------------------------------------------
$drugs.Cells.Item(2, 1).Value = 1
$drugs.Cells.Item(2, 2).Value = " استامینوفن"
$drugs.Cells.Item(2, 3).Value = " استامینوفن"
$drugs.Cells.Item(2, 4).Value = " Acetaminophen"
$drugs.Cells.Item(2, 5).Value = "تایلنول,پاراستامول"
$drugs.Cells.Item(2, 6).Value = " مسکن,تب‌بر"
$drugs.Cells.Item(2, 7).Value = " قرص,شربت"
$drugs.Cells.Item(2, 8).Value = " برای تسکین درد و کاهش تب استفاده می‌شود."
$drugs.Cells.Item(2, 9).Value = " تهوع,سردرد,خستگی"
$drugs.Cells.Item(2, 10).Value = " وارفارین,الکل"
$drugs.Cells.Item(2, 11).Value = " بیماری کبدی"
$drugs.Cells.Item(2, 12).Value = " در دمای اتاق و دور از نور مستقیم نگهداری شود."
$drugs.Cells.Item(2, 13).Value = "بزرگسالان: هر 4-6 ساعت 500-1000 میلی‌گرم"
$drugs.Cells.Item(2, 14).Value = "https://picsum.photos/seed/drug1/400/400"

# Column M (dosage) now holds longer text -> widen it, dropping the old
# auto "best fit" flag in favour of an explicit custom width.
$drugs.Columns.Item(13).ColumnWidth = 20.5

# --- New "Sheet1" worksheet, placed right after "drugs" --------------------
$newSheet = $wb.Sheets.Add([System.Type]::Missing, $drugs)

$newSheet.Cells.Item(1, 1).Value = 1
$newSheet.Cells.Item(2, 1).Value = " استامینوفن"
$newSheet.Cells.Item(3, 1).Value = " استامینوفن"
$newSheet.Cells.Item(4, 1).Value = " Acetaminophen"
$newSheet.Cells.Item(5, 1).Value = "تایلنول,پاراستامول"
$newSheet.Cells.Item(6, 1).Value = " مسکن,تب‌بر"
$newSheet.Cells.Item(7, 1).Value = " قرص,شربت"
$newSheet.Cells.Item(8, 1).Value = " برای تسکین درد و کاهش تب استفاده می‌شود."
$newSheet.Cells.Item(9, 1).Value = " تهوع,سردرد,خستگی"
$newSheet.Cells.Item(10, 1).Value = " وارفارین,الکل"
$newSheet.Cells.Item(11, 1).Value = " بیماری کبدی"
$newSheet.Cells.Item(12, 1).Value = " در دمای اتاق و دور از نور مستقیم نگهداری شود."
$newSheet.Cells.Item(13, 1).Value = "بزرگسالان: هر 4-6 ساعت 500-1000 میلی‌گرم"
$newSheet.Cells.Item(14, 1).Value = "https://picsum.photos/seed/drug1/400/400"

$newSheet.Range("A1:A14").Select() | Out-Null

# Leave "drugs" as the active sheet/selection, matching the new record row.
$drugs.Activate()
$drugs.Range("A2:N2").Select() | Out-Null
